$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price + 1h volume change) scraped on
# Sat Jun  3 07:51:55 UTC 2023. Rows 30/31 (Filecoin / InternetComputer)
# also swapped ranking order, so Coin name + Link are rewritten too.

$ws.Range('D2').Value = '27.180.69'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '1.904.25'
$ws.Range('E3').Value = '  +0.94%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '306.31'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5237'
$ws.Range('E7').Value = '  +1.58%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3766'
$ws.Range('E8').Value = '  +1.20%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07256'
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.15'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9000'
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08405'
$ws.Range('E12').Value = '  +10.27%  '
$ws.Range('D13').Value = '1.892.19'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '94.84'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.288'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008618'
$ws.Range('E17').Value = '  +1.30%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '14.59'
$ws.Range('E18').Value = '  +1.67%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').Value = '27.225.79'
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.067'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = '2.146.53'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.61'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.432'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '147.30'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.287'
$ws.Range('E26').Value = '  +5.14%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.752'
$ws.Range('E27').Value = '  -2.24%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.17'
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '114.90'
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.920'
$ws.Range('E30').Value = '  -1.24%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.815'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09269'
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.8095'
$ws.Range('E33').Value = '  +6.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05067'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.242'
$ws.Range('E35').Value = '  +4.21%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.967'
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.384'
$ws.Range('E37').Value = '  +3.46%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.609'
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5732'
$ws.Range('E39').Value = '  +1.79%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01989'
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.075'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.651'
$ws.Range('E42').Value = '  +1.13%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.990'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '117.32'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1513'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4853'
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.20'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.613'
$ws.Range('E49').Value = '  +1.62%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '37.44'
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '64.05'
$ws.Range('E51').Value = '  +0.53%  '
